$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Insert a new row above row 8 (pushes old row 8 "Upload" -> row 9, etc.)
$ws.Rows.Item(8).Insert()

# New row 8: "Force" label in column A, FALSE booleans in B:J,
# copying the formatting/style of the row above (row 7).
$ws.Range("A7:J7").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Force"
$ws.Range("B8:J8").Value = $false

# Re-freeze the panes one row lower (row 8 insert pushed the freeze
# boundary from row 9 to row 10) and restore the active selection.
$excel.ActiveWindow.FreezePanes = $false
[void]$ws.Range("A11").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A8:XFD8").Select()
